# "Fruta / hortaliza, semanal"
#
# The weekly refresh re-shuffles the per-market-day rows of the sheet: each
# row's Fecha / Volumen / Precio minimo / Precio maximo / Precio promedio
# ponderado / Precio $/Kg (columns D, J, K, L, M, P) ends up holding the
# values that another row used to hold. Everything else (Mercado ID,
# Mercado, Region, Codreg, Categoria, Variedad, Calidad, Unidad de
# comercializacion, Origen, Kg o Unidades, Clasificacion) is identical on
# every data row, so it is untouched.
#
# Snapshot the current values first (so writes never clobber a value that
# still needs to be read for a later row), then re-map them in one pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 45
$cols = 4, 10, 11, 12, 13, 16   # D, J, K, L, M, P

# Row `n` (key) takes the six values that row `n` held before the edit
# (value) in the destination.
$rowMap = @{
    2 = 34; 3 = 23; 4 = 3;   5 = 17; 6 = 11; 7 = 44; 8 = 32; 9 = 18; 10 = 15
    11 = 29; 12 = 21; 13 = 4; 14 = 12; 15 = 19; 16 = 40; 17 = 6; 18 = 35
    19 = 31; 20 = 45; 21 = 36; 22 = 7; 23 = 39; 24 = 37; 25 = 13; 26 = 10
    27 = 26; 28 = 41; 29 = 27; 30 = 33; 31 = 42; 32 = 30; 33 = 8; 34 = 43
    35 = 25; 36 = 38; 37 = 5; 38 = 2; 39 = 20; 40 = 9; 41 = 22; 42 = 14
    43 = 28; 44 = 24; 45 = 16
}

# 1) Snapshot every candidate cell's current value.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Write each row's new values from the mapped source row's snapshot.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $rowMap[$r]
    if ($src -eq $r) { continue }
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $snapshot["$src-$c"]
    }
}
